$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "kalori_2015_Maret"
$ws.Range("D1").Value = "kalori_2015_September"
$ws.Range("E1").Value = "kalori_2016_Maret"
$ws.Range("F1").Value = "kalori_2016_September"
$ws.Range("G1").Value = "kalori_2017_Maret"
$ws.Range("H1").Value = "kalori_2017_September"
$ws.Range("I1").Value = "kalori_2018_Maret"
$ws.Range("J1").Value = "kalori_2018_September"
$ws.Range("K1").Value = "kalori_2019_Maret"
$ws.Range("L1").Value = "kalori_2019_September"
$ws.Range("M1").Value = "kalori_2020_Maret"
$ws.Range("N1").Value = "kalori_2021_Maret"
$ws.Range("O1").Value = "kalori_2021_September"
$ws.Range("P1").Value = "kalori_2022_Maret"
$ws.Range("Q1").Value = "kalori_2022_September"
$ws.Range("R1").Value = "kalori_2023_Maret"
$ws.Range("S1").Value = "kalori_2024_Maret"
